$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43 ---
$ws.Range("C43").Value = 4
$ws.Range("E43").Value = 5
$ws.Range("G43").Value = 2
$ws.Range("I43").Value = 4
$ws.Range("J43").Value = 7
$ws.Range("K43").Value = 11
$ws.Range("L43").Value = 16
$ws.Range("M43").Value = 3
$ws.Range("N43").Value = 6
$ws.Range("O43").Value = 6
$ws.Range("P43").Value = 6

# --- Row 44 ---
$ws.Range("C44").Value = 6
$ws.Range("E44").Value = 10
$ws.Range("G44").Value = 2
$ws.Range("I44").Value = 4
$ws.Range("J44").Value = 3
$ws.Range("K44").Value = 4
$ws.Range("L44").Value = 5
$ws.Range("M44").Value = 5
$ws.Range("N44").Value = 5
$ws.Range("O44").Value = 5
$ws.Range("P44").Value = 5

# --- Row 45 ---
$ws.Range("C45").Value = 8
$ws.Range("E45").Value = 25
$ws.Range("G45").Value = 4
$ws.Range("I45").Value = 8
$ws.Range("J45").Value = 15
$ws.Range("K45").Value = 25
$ws.Range("L45").Value = 38
$ws.Range("M45").Value = 5
$ws.Range("N45").Value = 5
$ws.Range("O45").Value = 5
$ws.Range("P45").Value = 5
$ws.Range("I45").Style = "Bad"

# --- Row 46 ---
$ws.Range("C46").Value = 10
$ws.Range("E46").Value = 50
$ws.Range("G46").Value = 17
$ws.Range("I46").Value = 21
$ws.Range("J46").Value = 25
$ws.Range("K46").Value = 29
$ws.Range("L46").Value = 34
$ws.Range("M46").Value = 8
$ws.Range("N46").Value = 18
$ws.Range("O46").Value = 22
$ws.Range("P46").Value = 63
$ws.Range("G46:P46").Style = "Bad"

# --- Row 47 ---
$ws.Range("C47").Value = 15
$ws.Range("E47").Value = 150
$ws.Range("G47").Value = 26
$ws.Range("I47").Value = 31
$ws.Range("J47").Value = 42
$ws.Range("K47").Value = 58
$ws.Range("L47").Value = 79
$ws.Range("M47").Value = 23
$ws.Range("N47").Value = 29
$ws.Range("O47").Value = 80
$ws.Range("P47").Value = 102
$ws.Range("G47:P47").Style = "Bad"

# --- Row 48 ---
$ws.Range("C48").Value = 20
$ws.Range("E48").Value = 75
$ws.Range("G48").Value = 15
$ws.Range("I48").Value = 22
$ws.Range("J48").Value = 35
$ws.Range("K48").Value = 54
$ws.Range("L48").Value = 79
$ws.Range("M48").Value = 13
$ws.Range("N48").Value = 24
$ws.Range("O48").Value = 34
$ws.Range("P48").Value = 57
$ws.Range("G48").Style = "Bad"
$ws.Range("I48:N48").Style = "Bad"
$ws.Range("O48:P48").Style = "Good"

# --- Row 49 ---
$ws.Range("C49").Value = 25
$ws.Range("E49").Value = 100
$ws.Range("G49").Value = 20
$ws.Range("I49").Value = 25
$ws.Range("J49").Value = 34
$ws.Range("K49").Value = 47
$ws.Range("L49").Value = 64
$ws.Range("M49").Value = 8
$ws.Range("N49").Value = 11
$ws.Range("O49").Value = 15
$ws.Range("P49").Value = 31
$ws.Range("I49:O49").Style = "Bad"
$ws.Range("G49").Style = "Good"
$ws.Range("P49").Style = "Good"

# --- View state: move selection to I4 ---
$ws.Range("I4").Select() | Out-Null

# --- Page setup (explicit orientation triggers pageSetup element) ---
$ws.PageSetup.Orientation = 1
